$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 103.300135879668
$ws.Range("M2").Value = 1388.770259900942
$ws.Range("N2").Value = 1865.161350058632
$ws.Range("O2").Value = 1388.770209898442
$ws.Range("P2").Value = 2296.195837059123
$ws.Range("Q2").Value = 0.4799943480101776
$ws.Range("S2").Value = 3.098159155625434
$ws.Range("T2").Value = 0.5273835029111082
$ws.Range("U2").Value = 0.9599886960203552
$ws.Range("Z2").Value = 1410.117700967209
$ws.Range("AA2").Value = 255.643187439665
$ws.Range("AB2").Value = 0.6083805423608668
$ws.Range("AC2").Value = 1369.184195577029
$ws.Range("AD2").Value = 12.20517888828285
$ws.Range("AE2").Value = 0.2396586910868928

$ws.Range("C3").Value = 103.3338255422548
$ws.Range("M3").Value = 1388.743166782617
$ws.Range("N3").Value = 1535.203446008299
$ws.Range("O3").Value = 1388.743116780117
$ws.Range("P3").Value = 1950.734449187061
$ws.Range("Q3").Value = 0.4919741488262277
$ws.Range("S3").Value = 2.942939466986254
$ws.Range("T3").Value = 0.5453347087019914
$ws.Range("U3").Value = 0.9839482976524554
$ws.Range("Z3").Value = 1410.077438336539
$ws.Range("AA3").Value = 223.4168593657085
$ws.Range("AB3").Value = 0.6878781355861898
$ws.Range("AC3").Value = 1369.196039466222
$ws.Range("AD3").Value = 3.760501237345336
$ws.Range("AE3").Value = 0.2457162244265287

$ws.Range("C4").Value = 103.3069264252813
$ws.Range("M4").Value = 1388.708646950132
$ws.Range("N4").Value = 1646.273567690826
$ws.Range("O4").Value = 1388.708596947632
$ws.Range("P4").Value = 2156.275771687106
$ws.Range("Q4").Value = 0.5043402133436925
$ws.Range("S4").Value = 2.539133609290107
$ws.Range("T4").Value = 0.5593361778146114
$ws.Range("U4").Value = 1.008680426687385
$ws.Range("Z4").Value = 1410.029892626208
$ws.Range("AA4").Value = 240.4895425815223
$ws.Range("AB4").Value = 0.6510790748080652
$ws.Range("AC4").Value = 1370.314427822489
$ws.Range("AD4").Value = 34.70685387213599
$ws.Range("AE4").Value = 0.2519692921272606

$ws.Range("C5").Value = 103.3010317965518
$ws.Range("M5").Value = 1388.682090464793
$ws.Range("N5").Value = 2473.763813115619
$ws.Range("O5").Value = 1388.682040462293
$ws.Range("P5").Value = 3286.543689122479
$ws.Range("Q5").Value = 0.5132275679797801
$ws.Range("S5").Value = 3.342942663867335
$ws.Range("T5").Value = 0.5510927029344644
$ws.Range("U5").Value = 1.02645513595956
$ws.Range("Z5").Value = 1410.023726701798
$ws.Range("AA5").Value = 359.9721309540363
$ws.Range("AB5").Value = 0.6528278673564269
$ws.Range("AC5").Value = 1370.312458948726
$ws.Range("AD5").Value = 47.48534995287748
$ws.Range("AE5").Value = 0.2564643413273823

$ws.Range("C6").Value = 103.252475863462
$ws.Range("M6").Value = 1388.8395482197
$ws.Range("N6").Value = 1547.949272596378
$ws.Range("O6").Value = 1388.8394982172
$ws.Range("P6").Value = 1820.218237122434
$ws.Range("Q6").Value = 0.4605597013831995
$ws.Range("S6").Value = 2.139296855640777
$ws.Range("T6").Value = 0.5155992713955349
$ws.Range("U6").Value = 0.9211194027663989
$ws.Range("Z6").Value = 1410.145378400677
$ws.Range("AA6").Value = 203.1441747114904
$ws.Range("AB6").Value = 0.6080119621150195
$ws.Range("AC6").Value = 1371.159632636185
$ws.Range("AD6").Value = 13.4412817807938
$ws.Range("AE6").Value = 0.2297923557885428

$ws.Range("C7").Value = 103.3060452136017
$ws.Range("M7").Value = 1388.679911993933
$ws.Range("N7").Value = 1878.118281975277
$ws.Range("O7").Value = 1388.679861991433
$ws.Range("P7").Value = 2487.28923961154
$ws.Range("Q7").Value = 0.5128095786982665
$ws.Range("S7").Value = 3.131710173606903
$ws.Range("T7").Value = 0.5450616074977804
$ws.Range("U7").Value = 1.025619157396533
$ws.Range("Z7").Value = 1410.009265590518
$ws.Range("AA7").Value = 275.816500157199
$ws.Range("AB7").Value = 0.6728392920520275
$ws.Range("AC7").Value = 1369.179449175322
$ws.Range("AD7").Value = 10.78770316953782
$ws.Range("AE7").Value = 0.2562731797252413

$ws.Range("C8").Value = 103.3032998517817
$ws.Range("M8").Value = 1388.665964192076
$ws.Range("N8").Value = 1383.807595223978
$ws.Range("O8").Value = 1388.665914189576
$ws.Range("P8").Value = 1829.198093043822
$ws.Range("Q8").Value = 0.5145501038564645
$ws.Range("S8").Value = 2.617684887504798
$ws.Range("T8").Value = 0.5315474342961506
$ws.Range("U8").Value = 1.029100207712929
$ws.Range("Z8").Value = 1410.003405004466
$ws.Range("AA8").Value = 199.3622425849727
$ws.Range("AB8").Value = 0.7052610027123137
$ws.Range("AC8").Value = 1371.225091295636
$ws.Range("AD8").Value = 14.84828063850859
$ws.Range("AE8").Value = 0.2571747799808753

$ws.Range("C9").Value = 103.3077507680102
$ws.Range("M9").Value = 1388.649499691568
$ws.Range("N9").Value = 2189.886277459186
$ws.Range("O9").Value = 1388.649449689068
$ws.Range("P9").Value = 2925.083423797566
$ws.Range("Q9").Value = 0.5113938462770532
$ws.Range("S9").Value = 3.959809784928344
$ws.Range("T9").Value = 0.5738638351869245
$ws.Range("U9").Value = 1.022787692554106
$ws.Range("Z9").Value = 1410.009798128594
$ws.Range("AA9").Value = 342.6787436873864
$ws.Range("AB9").Value = 0.6649263291716057
$ws.Range("AC9").Value = 1371.781029057047
$ws.Range("AD9").Value = 11.12799435331445
$ws.Range("AE9").Value = 0.2555933995821887
